$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Terminal Hortofrutícola Agro Chillán - Berenjena)
# needs to be inserted as row 25, pushing the existing rows 25-46 down to 26-47.
$ws.Rows.Item(25).Insert()

$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44810
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112001
$ws.Range("G25").Value = "Berenjena"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 12000
$ws.Range("L25").Value = 13000
$ws.Range("M25").Value = 12500
$ws.Range("N25").Value = "$/caja 60 unidades"
$ws.Range("O25").Value = "Región de Arica y Parinacota"
$ws.Range("P25").Value = 208
$ws.Range("Q25").Value = 60
$ws.Range("R25").Value = "Hortaliza"

# Preserve the date number format used by the other rows in column D.
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
